# DP-610: LN_INTF generate test scripts
#
# Renames the two working sheets to their new "_INIT" / "_UPDATE" names,
# fixes up the per-sheet hidden _FilterDatabase defined names so they keep
# pointing at the right (renamed) sheet/range in the right order, updates
# the source-table reference used to build the "SELECT * FROM ..." query on
# the LN_INTF_INIT sheet (CYC_INTF -> LN_INTF), and moves the active-cell
# selection on the LN_INTF_UPDATE sheet.

$wb = $excel.ActiveWorkbook

$wsInit   = $wb.Worksheets.Item(1)   # was "LN_INTF"
$wsUpdate = $wb.Worksheets.Item(2)   # was "Case U_I_D"

# 1) Rename the sheets.
$wsInit.Name   = "LN_INTF_INIT"
$wsUpdate.Name = "LN_INTF_UPDATE"

# 2) Rebuild the hidden AutoFilter-database defined names so that
#    LN_INTF_INIT (localSheetId 0) is listed before LN_INTF_UPDATE
#    (localSheetId 1), each still hidden.
while ($wb.Names.Count -gt 0) {
  $wb.Names.Item(1).Delete()
}

$wsInit.Names.Add("_xlnm._FilterDatabase", "=LN_INTF_INIT!`$A`$1:`$U`$55") | Out-Null
$wsUpdate.Names.Add("_xlnm._FilterDatabase", "=LN_INTF_UPDATE!`$A`$1:`$U`$1") | Out-Null

$wb.Names.Item(1).Visible = $false
$wb.Names.Item(2).Visible = $false

# 3) The metadata-query helper table used CYC_INTF as a placeholder source
#    table name; point it at the real LN_INTF table (the dependent
#    "SELECT * FROM ..." formula in K2 recalculates automatically).
$wsInit.Range("O2").Value = "LN_INTF"

# 4) Move the selection on LN_INTF_UPDATE to H4, then restore LN_INTF_INIT
#    as the active sheet/tab (matches original tab selection).
$wsUpdate.Activate()
$wsUpdate.Range("H4").Select()
$wsInit.Activate()
